# Remove config file handling
# Appends a new row (row 35) to each of the four worksheets, duplicating the
# last existing row (row 34) but with the timestamp in column A advanced by
# one hour (17:xx:xx -> 18:xx:xx). All other columns (B-I) are identical to
# row 34 on the same sheet.

$wb = $excel.ActiveWorkbook

$newRows = @{
    1 = @{
        A = "2025-03-05 18:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    }
    2 = @{
        A = "2025-03-05 18:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    }
    3 = @{
        A = "2025-03-05 18:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    }
    4 = @{
        A = "2025-03-05 18:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $row = $newRows[$i]

    $ws.Cells.Item(35, 1).Value = $row.A
    $ws.Cells.Item(35, 2).Value = $row.B
    $ws.Cells.Item(35, 3).Value = $row.C
    $ws.Cells.Item(35, 4).Value = $row.D
    $ws.Cells.Item(35, 5).Value = $row.E
    $ws.Cells.Item(35, 6).Value = $row.F

    # Column G holds a 24-digit numeric string; writing it straight would be
    # silently coerced to a floating point number (losing precision), so
    # force text interpretation via the number format, then restore the
    # default "Normal" style so no stray formatting is left behind.
    $gCell = $ws.Cells.Item(35, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row.G
    $gCell.Style = "Normal"

    $ws.Cells.Item(35, 8).Value = $row.H
    $ws.Cells.Item(35, 9).Value = $row.I
}

Write-Output "Added row 35 to all sheets"
